# Fruta / hortaliza, semanal
# Update the weekly price rows (3-9, 11) with refreshed data.
# Row 10 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value2 = 44418
$ws.Range("L3").Value2 = "Especial"
$ws.Range("M3").Value2 = 100
$ws.Range("N3").Value2 = 8000
$ws.Range("O3").Value2 = 8000
$ws.Range("P3").Value2 = 8000
$ws.Range("Q3").Value2 = "$/caja 15 kilos granel"
$ws.Range("S3").Value2 = 533
$ws.Range("T3").Value2 = 15

# Row 4
$ws.Range("D4").Value2 = 44208
$ws.Range("L4").Value2 = "Especial"
$ws.Range("M4").Value2 = 70
$ws.Range("N4").Value2 = 24000
$ws.Range("O4").Value2 = 24000
$ws.Range("P4").Value2 = 24000
$ws.Range("Q4").Value2 = "$/caja 15 kilos granel"
$ws.Range("S4").Value2 = 1600
$ws.Range("T4").Value2 = 15

# Row 5
$ws.Range("D5").Value2 = 44411
$ws.Range("L5").Value2 = "Primera"
$ws.Range("M5").Value2 = 210
$ws.Range("N5").Value2 = 8000
$ws.Range("O5").Value2 = 8000
$ws.Range("P5").Value2 = 8000
$ws.Range("Q5").Value2 = "$/bandeja 8 kilos"
$ws.Range("S5").Value2 = 1000
$ws.Range("T5").Value2 = 8

# Row 6
$ws.Range("D6").Value2 = 44264
$ws.Range("L6").Value2 = "Calibre 100"
$ws.Range("M6").Value2 = 50
$ws.Range("N6").Value2 = 20000
$ws.Range("O6").Value2 = 20000
$ws.Range("P6").Value2 = 20000
$ws.Range("Q6").Value2 = "$/caja 18 kilos embalada"
$ws.Range("S6").Value2 = 1111
$ws.Range("T6").Value2 = 18

# Row 7
$ws.Range("D7").Value2 = 44217
$ws.Range("L7").Value2 = "Primera"
$ws.Range("M7").Value2 = 55
$ws.Range("N7").Value2 = 18000
$ws.Range("O7").Value2 = 18000
$ws.Range("P7").Value2 = 18000
$ws.Range("Q7").Value2 = "$/caja 18 kilos granel"
$ws.Range("S7").Value2 = 1000
$ws.Range("T7").Value2 = 18

# Row 8
$ws.Range("D8").Value2 = 44601
$ws.Range("L8").Value2 = "Primera"
$ws.Range("M8").Value2 = 30
$ws.Range("N8").Value2 = 28000
$ws.Range("O8").Value2 = 28000
$ws.Range("P8").Value2 = 28000
$ws.Range("Q8").Value2 = "$/caja 18 kilos granel"
$ws.Range("S8").Value2 = 1556
$ws.Range("T8").Value2 = 18

# Row 9
$ws.Range("D9").Value2 = 44392
$ws.Range("L9").Value2 = "Especial"
$ws.Range("M9").Value2 = 500
$ws.Range("N9").Value2 = 7000
$ws.Range("O9").Value2 = 7000
$ws.Range("P9").Value2 = 7000
$ws.Range("Q9").Value2 = "$/bandeja 8 kilos"
$ws.Range("S9").Value2 = 875
$ws.Range("T9").Value2 = 8

# Row 11
$ws.Range("D11").Value2 = 44427
$ws.Range("L11").Value2 = "Primera"
$ws.Range("M11").Value2 = 55
$ws.Range("N11").Value2 = 7000
$ws.Range("O11").Value2 = 7000
$ws.Range("P11").Value2 = 7000
$ws.Range("Q11").Value2 = "$/caja 15 kilos granel"
$ws.Range("S11").Value2 = 467
$ws.Range("T11").Value2 = 15
